# Ver-2.1.1
# 1. Add `num_of_components` metric as new column K.
# 2. Populate results for all data rows (2-93).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell K1 -----------------------------------------------------
$ws.Range("K1").Value = "num_of_components"

# Copy the header formatting (bold font, borders, centered alignment)
# from the adjacent J1 header cell so K1 matches the rest of the header row.
$ws.Range("J1").Copy()
$ws.Range("K1").PasteSpecial(-4122)

# --- Data cells K2:K93 ---------------------------------------------------
$numComponents = @(1,2,1,1,2,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,1,4,3,1,1,1,1,1,1,1,1,1,1,1,2,1,1,1,1,1,2,2,1,1,3,1,1,1,2,1,1,1,1,1,2,1,1,2,2,1,1,1,1,1,1,1,1,1,1,1,1,1,1,2,1,1,1,1,1,3,1,1,2,1,1,1,2,1,1,1,1)

$rowCount = $numComponents.Length
$arr = New-Object 'object[,]' $rowCount,1
for ($i = 0; $i -lt $rowCount; $i++) {
    $arr[$i,0] = $numComponents[$i]
}

$ws.Range("K2:K93").Value = $arr
